$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text representation
# (values like "1.00", "0.547", "51.420.50" must not be coerced to numbers)
$ws.Columns("D:E").NumberFormat = "@"

$ws.Range("D2").Value = "51.382.60"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "2.772.51"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "353.31"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").Value = "108.22"
$ws.Range("E6").Value = "  -1.22%  "
$ws.Range("D7").Value = "0.547"
$ws.Range("E7").Value = "  -3.05%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").Value = "39.65"
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("E11").Value = "  +3.03%  "
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").Value = "20.02"
$ws.Range("E12").Value = "  +3.44%  "
$ws.Range("B13").Value = "Dogecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D13").Value = "0.0833"
$ws.Range("E13").Value = "  -2.04%  "
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").Value = "3.210.16"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").Value = "2.776.67"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("D17").Value = "0.921"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").Value = "51.365.10"
$ws.Range("D19").Value = "7.60"
$ws.Range("E19").Value = "  +2.39%  "
$ws.Range("D20").Value = "3.10"
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("E22").Value = "  -1.41%  "
$ws.Range("D23").Value = "69.75"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "265.36"
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D27").Value = "25.97"
$ws.Range("E27").Value = "  -1.82%  "
$ws.Range("D28").Value = "0.161"
$ws.Range("E28").Value = "  +12.61%  "
$ws.Range("D29").Value = "10.19"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").Value = "36.25"
$ws.Range("E30").Value = "  +7.09%  "
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("E32").Value = "  +8.50%  "
$ws.Range("D33").Value = "51.90"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("D34").Value = "0.0448"
$ws.Range("E34").Value = "  -2.51%  "
$ws.Range("D35").Value = "5.51"
$ws.Range("E35").Value = "  +5.60%  "
$ws.Range("E36").Value = "  -2.14%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").Value = "18.17"
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("E39").Value = "  -2.34%  "
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("D43").Value = "120.61"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").Value = "21.93"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("E45").Value = "  -2.09%  "
$ws.Range("D46").Value = "2.096.07"
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("D48").Value = "2.31"
$ws.Range("E48").Value = "  +4.07%  "
$ws.Range("E49").Value = "  -4.54%  "
$ws.Range("D50").Value = "0.902"
$ws.Range("E50").Value = "  -2.15%  "
$ws.Range("E51").Value = "  +8.36%  "
